{"js": "// Change \"IdRegla = 39\" -> \"IdRegla = 70\" and \"IdRegla = 40\" -> \"IdRegla = 71\".\n// The target OOXML splits the original single bold run into two bold runs:\n// one with \"IdRegla = \" (trailing space preserved) and one with just the\n// new number. We rebuild that structure with insertOoxml so the run is\n// genuinely split (a plain insertText/\"Replace\" on a single range gets\n// re-coalesced into one run when the document is saved).\nfunction buildRunsOoxml(prefixText, numberText) {\n  return (\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" ' +\n    'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body><w:p>' +\n    '<w:r><w:rPr><w:b/></w:rPr><w:t xml:space=\"preserve\">' + prefixText + '</w:t></w:r>' +\n    '<w:r><w:rPr><w:b/></w:rPr><w:t>' + numberText + '</w:t></w:r>' +\n    '</w:p></w:body></w:document>' +\n    '</pkg:xmlData></pkg:part></pkg:package>'\n  );\n}\n\nasync function replaceIdRegla(oldText, newNumber) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    return;\n  }\n\n  results.items[0].insertOoxml(buildRunsOoxml(\"IdRegla = \", newNumber), \"Replace\");\n  await context.sync();\n}\n\nawait replaceIdRegla(\"IdRegla = 39\", \"70\");\nawait replaceIdRegla(\"IdRegla = 40\", \"71\");\n", "ps1": "# Change \"IdRegla = 39\" -> \"IdRegla = 70\" and \"IdRegla = 40\" -> \"IdRegla = 71\".\n# The target OOXML splits the original single bold run (e.g. \"IdRegla = 39\")\n# into two bold runs: \"IdRegla = \" (trailing space preserved) and the new\n# number (\"70\"/\"71\"). A plain Range.Text assignment re-coalesces into a\n# single run when saved, so we rebuild the paragraph via InsertXML, reusing\n# the paragraph's own existing attributes/pPr (read back from the live\n# document) so nothing else about the paragraph is altered.\n\n$d = $word.ActiveDocument\n\nfunction Replace-IdRegla($oldNumber, $newNumber) {\n  $range = $d.Content\n  $find = $range.Find\n  $find.Text = \"IdRegla = \" + $oldNumber\n  $found = $find.Execute()\n  if (-not $found) {\n    return\n  }\n\n  # Read back the enclosing paragraph's real opening-tag attributes and\n  # pPr block so the rebuilt paragraph keeps them unchanged.\n  $para = $range.Paragraphs.Item(1)\n  $pxml = $para.Range.WordOpenXML\n\n  $pOpenStart = $pxml.IndexOf(\"<w:p \")\n  if ($pOpenStart -lt 0) {\n    $pOpenStart = $pxml.IndexOf(\"<w:p>\")\n    $pOpenAttrs = \"\"\n    $pOpenEnd = $pOpenStart + 4\n  } else {\n    $pOpenEnd = $pxml.IndexOf(\">\", $pOpenStart)\n    $pOpenAttrs = $pxml.Substring($pOpenStart + 4, $pOpenEnd - $pOpenStart - 4)\n    # Drop w14:paraId/w14:textId (computed, not present in the stored file).\n    $pOpenAttrs = $pOpenAttrs -replace 'w14:\\w+=\"[^\"]*\"\\s*', ''\n    $pOpenAttrs = $pOpenAttrs.Trim()\n  }\n\n  $pPrBlock = \"\"\n  $pprStart = $pxml.IndexOf(\"<w:pPr>\", $pOpenEnd)\n  if ($pprStart -ge $pOpenEnd -and $pprStart -le ($pOpenEnd + 10)) {\n    $pprEnd = $pxml.IndexOf(\"</w:pPr>\", $pprStart) + \"</w:pPr>\".Length\n    $pPrBlock = $pxml.Substring($pprStart, $pprEnd - $pprStart)\n  }\n\n  $pOpenTag = \"<w:p\"\n  if ($pOpenAttrs -ne \"\") {\n    $pOpenTag = $pOpenTag + \" \" + $pOpenAttrs\n  }\n  $pOpenTag = $pOpenTag + \">\"\n\n  $xml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" ' +\n    'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body>' + $pOpenTag + $pPrBlock +\n    '<w:r><w:rPr><w:b/></w:rPr><w:t xml:space=\"preserve\">IdRegla = </w:t></w:r>' +\n    '<w:r><w:rPr><w:b/></w:rPr><w:t>' + $newNumber + '</w:t></w:r>' +\n    '</w:p></w:body></w:document>' +\n    '</pkg:xmlData></pkg:part></pkg:package>'\n\n  $range.InsertXML($xml)\n}\n\nReplace-IdRegla \"39\" \"70\"\nReplace-IdRegla \"40\" \"71\"\n"}
